# Auto-generated edit script applying scheduled market-data refresh to Chocobo_Profits sheets.
# For each affected Leve row, columns H-N (currentAveragePrice.., LevePriceNQ/HQ, LeveProfitNQ/HQ)
# are updated to the latest computed values. Some cells are cleared entirely (no profit computed)
# or newly populated, matching the upstream data pull.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 282.66666
$ws.Range("I33").Value = 185.93333
$ws.Range("J33").Value = 766.3333
$ws.Range("K33").Value = 185.93333
$ws.Range("L33").Value = 766.3333
$ws.Range("M33").Value = 43.06666999999999
$ws.Range("N33").Value = -1224.3333

$ws.Range("H40").Value = 1639.8
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 1639.8
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 1639.8
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -1989.8

$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("M125").ClearContents()
$ws.Range("N125").ClearContents()

$ws.Range("H136").Value = 51168
$ws.Range("J136").Value = 51168
$ws.Range("L136").Value = 51168
$ws.Range("N136").Value = -61368

$ws.Range("H137").Value = 747116.4
$ws.Range("I137").Value = 2073979
$ws.Range("J137").Value = 2778.756
$ws.Range("K137").Value = 6221937
$ws.Range("L137").Value = 8336.268
$ws.Range("M137").Value = -6219387
$ws.Range("N137").Value = -13436.268

$ws.Range("H138").Value = 1880.5238
$ws.Range("I138").Value = 769.5714
$ws.Range("J138").Value = 4102.4287
$ws.Range("K138").Value = 2308.7142
$ws.Range("L138").Value = 12307.2861
$ws.Range("M138").Value = 2831.2858
$ws.Range("N138").Value = -22587.2861


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 38000
$ws.Range("J23").Value = 38000
$ws.Range("L23").Value = 38000
$ws.Range("N23").Value = -38518

$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()

$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()

$ws.Range("H113").Value = 49980
$ws.Range("J113").Value = 49980
$ws.Range("L113").Value = 49980
$ws.Range("N113").Value = -58658


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 5786539
$ws.Range("I7").Value = 7009201
$ws.Range("J7").Value = 5022375
$ws.Range("K7").Value = 7009201
$ws.Range("L7").Value = 5022375
$ws.Range("M7").Value = -7009088
$ws.Range("N7").Value = -5022601

$ws.Range("H63").Value = 37000
$ws.Range("J63").Value = 37000
$ws.Range("L63").Value = 37000
$ws.Range("N63").Value = -38372

$ws.Range("H66").Value = 37000
$ws.Range("J66").Value = 37000
$ws.Range("L66").Value = 111000
$ws.Range("N66").Value = -117864

$ws.Range("H94").Value = 843.26666
$ws.Range("I94").Value = 715.36365
$ws.Range("K94").Value = 715.36365
$ws.Range("M94").Value = -264.36365

$ws.Range("H112").Value = 29970
$ws.Range("J112").Value = 29970
$ws.Range("L112").Value = 29970
$ws.Range("N112").Value = -32924


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H47").Value = 55000
$ws.Range("J47").Value = 55000
$ws.Range("L47").Value = 55000
$ws.Range("N47").Value = -56132

$ws.Range("H68").Value = 75995.39999999999
$ws.Range("J68").Value = 75995.39999999999
$ws.Range("L68").Value = 75995.39999999999
$ws.Range("N68").Value = -77493.39999999999

$ws.Range("H71").Value = 75995.39999999999
$ws.Range("J71").Value = 75995.39999999999
$ws.Range("L71").Value = 227986.2
$ws.Range("N71").Value = -235474.2


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6825.9287
$ws.Range("I70").Value = 6175.1514
$ws.Range("J70").Value = 9212.111000000001
$ws.Range("K70").Value = 6175.1514
$ws.Range("L70").Value = 9212.111000000001
$ws.Range("M70").Value = -5905.1514
$ws.Range("N70").Value = -9752.111000000001

$ws.Range("H73").Value = 6825.9287
$ws.Range("I73").Value = 6175.1514
$ws.Range("J73").Value = 9212.111000000001
$ws.Range("K73").Value = 6175.1514
$ws.Range("L73").Value = 9212.111000000001
$ws.Range("M73").Value = -5239.1514
$ws.Range("N73").Value = -11084.111

$ws.Range("H75").Value = 34000
$ws.Range("J75").Value = 34000
$ws.Range("L75").Value = 34000
$ws.Range("N75").Value = -35748

$ws.Range("H78").Value = 34000
$ws.Range("J78").Value = 34000
$ws.Range("L78").Value = 102000
$ws.Range("N78").Value = -110736

$ws.Range("H93").Value = 9799.6
$ws.Range("J93").Value = 9799.6
$ws.Range("L93").Value = 9799.6
$ws.Range("N93").Value = -13543.6

$ws.Range("H113").Value = 2013.875
$ws.Range("I113").Value = 2085.1667
$ws.Range("K113").Value = 2085.1667
$ws.Range("M113").Value = 84.83329999999978

$ws.Range("H119").Value = 34900
$ws.Range("J119").Value = 34900
$ws.Range("L119").Value = 34900
$ws.Range("N119").Value = -44576


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7333.3335
$ws.Range("I7").Value = 2000
$ws.Range("J7").Value = 10000
$ws.Range("K7").Value = 2000
$ws.Range("L7").Value = 10000
$ws.Range("M7").Value = -1888
$ws.Range("N7").Value = -10224

$ws.Range("H74").Value = 43776.285
$ws.Range("J74").Value = 43776.285
$ws.Range("L74").Value = 43776.285
$ws.Range("N74").Value = -45772.285

$ws.Range("H77").Value = 43776.285
$ws.Range("J77").Value = 43776.285
$ws.Range("L77").Value = 131328.855
$ws.Range("N77").Value = -141312.855

$ws.Range("H93").Value = 5052528
$ws.Range("I93").Value = 5849716
$ws.Range("K93").Value = 5849716
$ws.Range("M93").Value = -5848468

$ws.Range("H122").Value = 6316.6665
$ws.Range("I122").Value = 5114.2856
$ws.Range("J122").Value = 8000
$ws.Range("K122").Value = 15342.8568
$ws.Range("L122").Value = 24000
$ws.Range("M122").Value = -12892.8568
$ws.Range("N122").Value = -28900

$ws.Range("H126").Value = 7333.3335
$ws.Range("I126").Value = 2000
$ws.Range("J126").Value = 10000
$ws.Range("K126").Value = 6000
$ws.Range("L126").Value = 30000
$ws.Range("M126").Value = -3530
$ws.Range("N126").Value = -34940


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 42500
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 42500
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 42500
$ws.Range("M75").ClearContents()
$ws.Range("N75").Value = -44372

$ws.Range("H78").Value = 42500
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 42500
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 127500
$ws.Range("N78").Value = -136860

